# Update UBID values in column E (shared strings) per "update ubids in example files"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "Z01TDR2Z+7ES-Z01TDR2Z+HX7-Z01TDR2Z+UAX"
$ws.Range("E3").Value = "M7RZ35FK+6LL-M7RZ35FK+9GS-M7RZ35FK+H0V"
$ws.Range("E4").Value = "1TVB419R+QBO-1TVB419R+CFW-1TVB419R+R9S"
$ws.Range("E5").Value = "ELPXXGCE+X4C-ELPXXGCE+3ZR-ELPXXGCE+78L"
$ws.Range("E6").Value = "BKSG43YH+3HW-BKSG43YH+WCL-BKSG43YH+DDD"
$ws.Range("E7").Value = "VEYQCG4R+Q5P-VEYQCG4R+X18-VEYQCG4R+A8U"
$ws.Range("E8").Value = "2YWG8HDH+H0D-2YWG8HDH+0U8-2YWG8HDH+IGU"
$ws.Range("E9").Value = "WW2YKUX2+FVE-WW2YKUX2+8SH-WW2YKUX2+3K2"
$ws.Range("E10").Value = "VQADDOC3+V8E-VQADDOC3+XBF-VQADDOC3+2EE"
$ws.Range("E11").Value = "ZIUC82DT+4X5-ZIUC82DT+C4M-ZIUC82DT+YK9"
$ws.Range("E12").Value = "4RFKPAAM+R47-4RFKPAAM+TTP-4RFKPAAM+19B"
$ws.Range("E13").Value = "LUFXFMVJ+8XY-LUFXFMVJ+QJG-LUFXFMVJ+QQ1"
$ws.Range("E14").Value = "VZTP3DRR+K1J-VZTP3DRR+X4E-VZTP3DRR+E8V"
$ws.Range("E15").Value = "QCILPX9G+1MT-QCILPX9G+NPX-QCILPX9G+HNH"

# Re-select the UBID column range, matching the selection left after the edit
[void]$ws.Range("E2:E15").Select()

# Re-fit column widths to the new content (best effort; engine rounds to 1/6-character steps)
$ws.Columns.Item(1).ColumnWidth = 19.166666666666668
$ws.Columns.Item(2).ColumnWidth = 7.498697916666667
$ws.Columns.Item(3).ColumnWidth = 9.330729166666666
$ws.Columns.Item(4).ColumnWidth = 13.166666666666666
$ws.Columns.Item(5).ColumnWidth = 40.498697916666664
$ws.Columns.Item(6).ColumnWidth = 17.498697916666668
$ws.Columns.Item(7).ColumnWidth = 17.330729166666668
$ws.Columns.Item(8).ColumnWidth = 3.4986979166666665
$ws.Columns.Item(9).ColumnWidth = 28.998697916666668
$ws.Columns.Item(10).ColumnWidth = 24.666666666666668
$ws.Columns.Item(11).ColumnWidth = 10.330729166666666
$ws.Columns.Item(12).ColumnWidth = 10.330729166666666
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667
$ws.Columns.Item(14).ColumnWidth = 10.330729166666666
$ws.Columns.Item(15).ColumnWidth = 12.830729166666666
$ws.Columns.Item(16).ColumnWidth = 14.998697916666666
$ws.Columns.Item(17).ColumnWidth = 11.830729166666666
$ws.Columns.Item(18).ColumnWidth = 38.330729166666664
$ws.Columns.Item(19).ColumnWidth = 11.830729166666666
$ws.Columns.Item(20).ColumnWidth = 10.166666666666666
$ws.Columns.Item(21).ColumnWidth = 9.998697916666666
$ws.Columns.Item(22).ColumnWidth = 7.666666666666667
$ws.Columns.Item(23).ColumnWidth = 7.666666666666667
